# RAG+C Q&A workbook — add an "answer_edited" column.
#
# The author duplicated the "answer" column (C) into a brand-new column D
# ("answer_edited"), pushing the old D/E/F (response_no_context,
# response_context, source) one slot to the right -> E/F/G. Two rows got a
# genuinely edited answer text (row 24 -> "Three", row 44 -> "Thirty
# days."); every other row's new D cell is just a copy of its C cell. The
# two stray hyperlinks that lived on the old "source" column were also
# removed, and the view was left scrolled down near the bottom of the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column C ("answer") into a freshly inserted column D
# ("answer_edited"); this shifts the old D/E/F -> E/F/G and copies both
# values and formatting along the way.
$ws.Columns("C:C").Copy()
$ws.Columns("D:D").Insert()

# Column D needs to be as wide as B/C (the paste above leaves it at the
# sheet's default column width).
$ws.Columns("D:D").ColumnWidth = 47.1640625

# Header + the two rows whose edited answer differs from the original.
$ws.Range("D1").Value = "answer_edited"
$ws.Range("D24").Value = "Three"
$ws.Range("D44").Value = "Thirty days."

# Drop the old hyperlinks entirely (previously anchored on the "source"
# column).
$ws.Hyperlinks.Delete()

# Leave the view parked where the editor last left it (scrolled down near
# the bottom of the sheet, one column over).
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D52").Select()

Write-Output "done"
